$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (index 1): plain value updates, no structural change.
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Cells.Item(2,6).Value = 292
$wsExpo.Cells.Item(3,6).Value = 19
$wsExpo.Cells.Item(4,6).Value = 7845
$wsExpo.Cells.Item(5,6).Value = 5734
$wsExpo.Cells.Item(6,6).Value = 476
$wsExpo.Cells.Item(8,6).Value = 12
$wsExpo.Cells.Item(10,6).Value = 265
$wsExpo.Cells.Item(11,6).Value = 295
$wsExpo.Cells.Item(12,6).Value = 61

# ---------------------------------------------------------------------
# Sheet "演出" (index 2): insert a new event row before the existing
# "2024-11-09" row (old row 3), pushing the rest down by one.
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Rows.Item(3).Insert()
# Match the sequence-number column formatting used by the other data rows.
$wsShow.Cells.Item(2,1).Copy()
$wsShow.Cells.Item(3,1).PasteSpecial(-4122)

$wsShow.Cells.Item(3,1).Value = 2
$wsShow.Cells.Item(3,2).Value = "'2024-11-08"
$wsShow.Cells.Item(3,3).Value = "合肥·松井祐贵 2024《阳光之旅》指弹吉他音乐会"
$wsShow.Cells.Item(3,4).Value = "宁国南路与水阳江路交口罍街二期15号楼安徽原创音乐基地3楼 OTW LIVEHOUSE"
$wsShow.Cells.Item(3,5).Value = "2024.11.08 19:30-11.08 21:00"
$wsShow.Cells.Item(3,6).Value = 0
$wsShow.Cells.Item(3,7).Value = 220
$wsShow.Cells.Item(3,8).Value = "https://show.bilibili.com/platform/detail.html?id=92768"
$wsShow.Cells.Item(3,9).Value = "//i1.hdslb.com/bfs/openplatform/202409/OU2qWxgM1727082424391.jpeg"

# Renumber the sequence column for the rows that shifted down.
$wsShow.Cells.Item(4,1).Value = 3
$wsShow.Cells.Item(5,1).Value = 4

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4): same value updates as 展览 for the first
# rows, plus the same new row inserted before the "2024-11-09" row
# (old row 12), pushing the following rows down by one.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Cells.Item(2,6).Value = 292
$wsAll.Cells.Item(3,6).Value = 19
$wsAll.Cells.Item(4,6).Value = 7845
$wsAll.Cells.Item(5,6).Value = 5734
$wsAll.Cells.Item(6,6).Value = 476
$wsAll.Cells.Item(8,6).Value = 12
$wsAll.Cells.Item(10,6).Value = 265

$wsAll.Rows.Item(12).Insert()
$wsAll.Cells.Item(2,1).Copy()
$wsAll.Cells.Item(12,1).PasteSpecial(-4122)

$wsAll.Cells.Item(12,1).Value = 10
$wsAll.Cells.Item(12,2).Value = "'2024-11-08"
$wsAll.Cells.Item(12,3).Value = "合肥·松井祐贵 2024《阳光之旅》指弹吉他音乐会"
$wsAll.Cells.Item(12,4).Value = "宁国南路与水阳江路交口罍街二期15号楼安徽原创音乐基地3楼 OTW LIVEHOUSE"
$wsAll.Cells.Item(12,5).Value = "2024.11.08 19:30-11.08 21:00"
$wsAll.Cells.Item(12,6).Value = 0
$wsAll.Cells.Item(12,7).Value = 220
$wsAll.Cells.Item(12,8).Value = "https://show.bilibili.com/platform/detail.html?id=92768"
$wsAll.Cells.Item(12,9).Value = "//i1.hdslb.com/bfs/openplatform/202409/OU2qWxgM1727082424391.jpeg"

# Renumber the sequence column and apply the value bumps for rows that
# shifted down from 12-15 to 13-16.
$wsAll.Cells.Item(13,1).Value = 11
$wsAll.Cells.Item(14,1).Value = 12
$wsAll.Cells.Item(14,6).Value = 295
$wsAll.Cells.Item(15,1).Value = 13
$wsAll.Cells.Item(15,6).Value = 61
$wsAll.Cells.Item(16,1).Value = 14
